# Checkout slot execution changes
# Applies updates to the TC_Checkout sheet: a handful of previously-failing
# steps now pass (Actual/Status/Comment updated), the sign-up email input
# was changed to a real mailinator address, and the order-instructions
# element locator/value was corrected to "Delivery Instructions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_Checkout")
$ws.Activate()

# --- Row 6: SetText correct email id ---
$ws.Range("G6").Value = "userchandna7049580@mailinator.com"
$ws.Range("L6").Value = "SetText: userchandna7049580@mailinator.com"
$ws.Rows.Item(6).RowHeight = 75

# --- Row 17: Add product now passes ---
$ws.Range("M17").Value = "Pass"

# --- Row 38: Click Proceed to payment now passes ---
$ws.Range("L38").Value = "Click: null"
$ws.Range("M38").Value = "Pass"
$ws.Range("N38").Value = "-"

# --- Row 39: VerifyText delivery slot message now passes ---
$ws.Range("L39").Value = "VerifyText: Please select delivery slot."
$ws.Range("M39").Value = "Pass"
$ws.Range("N39").Value = "-"

# --- Row 44: Enter Order/Delivery instruction now passes ---
$ws.Range("F44").Value = "//textarea[normalize-space(@placeholder) = 'Delivery Instructions']"
$ws.Range("L44").Value = "SetText: user chandna"
$ws.Range("M44").Value = "Pass"
$ws.Range("N44").Value = "-"

# --- Row 51: Click Place Order button now passes ---
$ws.Range("L51").Value = "Click: null"
$ws.Range("M51").Value = "Pass"
$ws.Range("N51").Value = "-"

# --- Row 52: VerifyText card security code message now passes ---
$ws.Range("L52").Value = "VerifyText: Your card's security code is incomplete"
$ws.Range("M52").Value = "Pass"
$ws.Range("N52").Value = "-"

# --- Update the saved view/selection state (scrolled down one row, new active cell) ---
$ws.Range("A40").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G44").Select()
